# chore: update Sheets via scheduled runner
# Refreshes cached market-board price/profit figures on the per-job
# "Yojimbo Profits" leve tracker sheets (ALC, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 11057.2
$ws.Range("I43").Value = 13259
$ws.Range("J43").Value = 2250
$ws.Range("K43").Value = 13259
$ws.Range("L43").Value = 2250
$ws.Range("M43").Value = -13190
$ws.Range("N43").Value = -2388
$ws.Range("H111").Value = 1175
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 1175
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 3525
$ws.Range("M111").ClearContents()   # profit now blank (merged into L/N)
$ws.Range("N111").Value = -9659
$ws.Range("H127").Value = 1102.4445
$ws.Range("J127").Value = 1442.4546
$ws.Range("L127").Value = 4327.3638
$ws.Range("N127").Value = -14247.3638
$ws.Range("H137").Value = 2632.547
$ws.Range("I137").Value = 2640.658
$ws.Range("J137").Value = 2612
$ws.Range("K137").Value = 7921.974
$ws.Range("L137").Value = 7836
$ws.Range("M137").Value = -5371.974
$ws.Range("N137").Value = -12936
$ws.Range("H138").Value = 2500.9702
$ws.Range("I138").Value = 1243.1111
$ws.Range("J138").Value = 3961.7097
$ws.Range("K138").Value = 3729.3333
$ws.Range("L138").Value = 11885.1291
$ws.Range("M138").Value = 1410.6667
$ws.Range("N138").Value = -22165.1291
$ws.Range("H141").Value = 2315.2258
$ws.Range("I141").Value = 2379.1304
$ws.Range("J141").Value = 2131.5
$ws.Range("K141").Value = 7137.3912
$ws.Range("L141").Value = 6394.5
$ws.Range("M141").Value = -1957.3912
$ws.Range("N141").Value = -16754.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30294.695
$ws.Range("I31").Value = 46703.184
$ws.Range("J31").Value = 4509.9287
$ws.Range("K31").Value = 46703.184
$ws.Range("L31").Value = 4509.9287
$ws.Range("M31").Value = -46408.184
$ws.Range("N31").Value = -5099.9287
$ws.Range("H34").Value = 30294.695
$ws.Range("I34").Value = 46703.184
$ws.Range("J34").Value = 4509.9287
$ws.Range("K34").Value = 46703.184
$ws.Range("L34").Value = 4509.9287
$ws.Range("M34").Value = -46501.184
$ws.Range("N34").Value = -4913.9287
$ws.Range("H132").Value = 1362.4
$ws.Range("I132").Value = 967.24
$ws.Range("J132").Value = 3338.2
$ws.Range("K132").Value = 2901.72
$ws.Range("L132").Value = 10014.6
$ws.Range("M132").Value = -371.7200000000003
$ws.Range("N132").Value = -15074.6
$ws.Range("H134").Value = 999.02563
$ws.Range("I134").Value = 904.5
$ws.Range("K134").Value = 2713.5
$ws.Range("M134").Value = -178.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1333.56
$ws.Range("I5").Value = 312
$ws.Range("J5").Value = 2276.5386
$ws.Range("K5").Value = 936
$ws.Range("L5").Value = 6829.6158
$ws.Range("M5").Value = -824
$ws.Range("N5").Value = -7053.6158
$ws.Range("H34").Value = 1064.7241
$ws.Range("I34").Value = 959.875
$ws.Range("J34").Value = 1104.6666
$ws.Range("K34").Value = 2879.625
$ws.Range("L34").Value = 3313.9998
$ws.Range("M34").Value = -2795.625
$ws.Range("N34").Value = -3481.9998
$ws.Range("H39").Value = 1200.3636
$ws.Range("J39").Value = 1200.3636
$ws.Range("L39").Value = 3601.0908
$ws.Range("N39").Value = -4189.0908
$ws.Range("H122").Value = 1451.3549
$ws.Range("I122").Value = 1195.5714
$ws.Range("J122").Value = 1662
$ws.Range("K122").Value = 10760.1426
$ws.Range("L122").Value = 14958
$ws.Range("M122").Value = -8310.142600000001
$ws.Range("N122").Value = -19858
$ws.Range("H135").Value = 1333.56
$ws.Range("I135").Value = 312
$ws.Range("J135").Value = 2276.5386
$ws.Range("K135").Value = 2808
$ws.Range("L135").Value = 20488.8474
$ws.Range("M135").Value = -273
$ws.Range("N135").Value = -25558.8474

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1078
$ws.Range("I113").Value = 925.7143
$ws.Range("J113").Value = 1433.3334
$ws.Range("K113").Value = 925.7143
$ws.Range("L113").Value = 1433.3334
$ws.Range("M113").Value = 1244.2857
$ws.Range("N113").Value = -5773.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 267.33334
$ws.Range("H27").Value = 267.33334
$ws.Range("H40").Value = 3551.2646
$ws.Range("I40").Value = 3055.9583
$ws.Range("J40").Value = 4740
$ws.Range("K40").Value = 3055.9583
$ws.Range("L40").Value = 4740
$ws.Range("M40").Value = -2919.9583
$ws.Range("N40").Value = -5012
$ws.Range("H46").Value = 500.25
$ws.Range("I46").Value = 500.25
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 500.25
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -312.25    # new HQ-profit figure now populated
$ws.Range("N46").ClearContents()    # superseded by M46
$ws.Range("H94").Value = 29000
$ws.Range("J94").Value = 29000
$ws.Range("L94").Value = 29000
$ws.Range("N94").Value = -30352     # new HQ-profit figure now populated
$ws.Range("H122").Value = 2899.75
$ws.Range("I122").Value = 3044.5293
$ws.Range("J122").Value = 2735.6667
$ws.Range("K122").Value = 9133.5879
$ws.Range("L122").Value = 8207.000100000001
$ws.Range("M122").Value = -6683.5879
$ws.Range("N122").Value = -13107.0001
$ws.Range("H132").Value = 11179.071
$ws.Range("I132").Value = 15413
$ws.Range("J132").Value = 5533.8335
$ws.Range("K132").Value = 46239
$ws.Range("L132").Value = 16601.5005
$ws.Range("M132").Value = -43709
$ws.Range("N132").Value = -21661.5005
$ws.Range("H134").Value = 5000
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()   # superseded by existing M134
$ws.Range("H136").Value = 3494.2124
$ws.Range("I136").Value = 1749.5818
$ws.Range("J136").Value = 7332.4
$ws.Range("K136").Value = 5248.7454
$ws.Range("L136").Value = 21997.2
$ws.Range("M136").Value = -2698.7454
$ws.Range("N136").Value = -27097.2
$ws.Range("H137").Value = 42960
$ws.Range("I137").Value = 15000
$ws.Range("J137").Value = 49950
$ws.Range("K137").Value = 15000
$ws.Range("L137").Value = 49950
$ws.Range("M137").Value = -9900
$ws.Range("N137").Value = -60150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 54975
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()    # superseded, N55 already has the figure
$ws.Range("H126").Value = 286477.6
$ws.Range("I126").Value = 370955.6
$ws.Range("J126").Value = 1364.375
$ws.Range("K126").Value = 1112866.8
$ws.Range("L126").Value = 4093.125
$ws.Range("M126").Value = -1110396.8
$ws.Range("N126").Value = -9033.125
$ws.Range("H132").Value = 4966.1665
$ws.Range("I132").Value = 5476.815
$ws.Range("K132").Value = 16430.445
$ws.Range("M132").Value = -13900.445
$ws.Range("H136").Value = 709.28
$ws.Range("J136").Value = 839.375
$ws.Range("L136").Value = 2518.125
$ws.Range("N136").Value = -7618.125

